$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of accelerometer data to insert right after the header row (row 1),
# pushing the existing data down by 3 rows.
$topRows = @(
    @(1.45181941986084, -10.64814472198486, 5.583051204681396),
    @(12.47824478149414, -36.90848541259766, 23.18166542053223),
    @(4.87645149230957, -11.23874282836914, 11.8306131362915)
)

# Insert 3 new blank rows starting at row 2 (shifts current rows 2..21 down to 5..24)
$insertRange = $ws.Range("A2:C4")
$insertRange.Insert()
$insertRange.ClearFormats()

for ($i = 0; $i -lt $topRows.Count; $i++) {
    $r = 2 + $i
    $ws.Cells.Item($r, 1).Value = $topRows[$i][0]
    $ws.Cells.Item($r, 2).Value = $topRows[$i][1]
    $ws.Cells.Item($r, 3).Value = $topRows[$i][2]
}

# New rows of accelerometer data appended after the (now shifted) last row.
# Original data ended at row 21, now at row 24, so new rows start at 25.
$bottomRows = @(
    @(-8.727315902709961, -7.33059024810791, 0.0287117958068847),
    @(-17.04781913757324, -51.86902236938477, 16.61567497253418),
    @(27.31503295898437, 1.940977096557617, 14.65385055541992),
    @(3.707320690155029, -6.758492469787598, 5.369882583618164),
    @(18.9058609008789, -66.31611633300781, 27.15024948120117),
    @(-72.50922393798828, 29.85196113586425, -13.03144264221191),
    @(-21.97713661193848, -12.47497940063477, 15.31039047241211)
)

$startRow = 25
for ($i = 0; $i -lt $bottomRows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $bottomRows[$i][0]
    $ws.Cells.Item($r, 2).Value = $bottomRows[$i][1]
    $ws.Cells.Item($r, 3).Value = $bottomRows[$i][2]
}
